$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.675.41"
$ws.Range("E2").Value = "  -4.05%  "

# Row 3
$ws.Range("D3").Value = "2.974.38"
$ws.Range("E3").Value = "  -5.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.27"
$ws.Range("E5").Value = "  -4.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.59"
$ws.Range("E6").Value = "  -5.91%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  +1.31%  "

# Row 9
$ws.Range("D9").Value = "2.985.13"
$ws.Range("E9").Value = "  -5.33%  "

# Row 10
$ws.Range("E10").Value = "  -2.26%  "

# Row 11
$ws.Range("E11").Value = "  -6.31%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  -2.53%  "

# Row 13
$ws.Range("D13").Value = "3.491.56"
$ws.Range("E13").Value = "  -5.30%  "

# Row 14
$ws.Range("E14").Value = "  -2.69%  "

# Row 15
$ws.Range("D15").Value = "61.738.33"
$ws.Range("E15").Value = "  -4.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.80"
$ws.Range("E16").Value = "  -4.58%  "

# Row 17
$ws.Range("D17").Value = "2.983.11"
$ws.Range("E17").Value = "  -5.24%  "

# Row 18
$ws.Range("E18").Value = "  -4.21%  "

# Row 19
$ws.Range("E19").Value = "  -0.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  -3.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.86"
$ws.Range("E21").Value = "  -4.96%  "

# Row 22
$ws.Range("E22").Value = "  -5.66%  "

# Row 23
$ws.Range("E23").Value = "  +0.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.97"
$ws.Range("E24").Value = "  -3.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.472"
$ws.Range("E25").Value = "  -2.11%  "

# Row 26
$ws.Range("D26").Value = "3.094.66"
$ws.Range("E26").Value = "  -5.55%  "

# Row 27
$ws.Range("E27").Value = "  -3.10%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0946"
$ws.Range("E28").Value = "  -6.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.34"
$ws.Range("E30").Value = "  -4.97%  "

# Row 31
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("E32").Value = "  -4.13%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.52"
$ws.Range("E33").Value = "  -3.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "161.17"
$ws.Range("E34").Value = "  +1.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.69"
$ws.Range("E35").Value = "  -2.59%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("E36").Value = "  -4.28%  "

# Row 37
$ws.Range("E37").Value = "  -2.34%  "

# Row 38
$ws.Range("E38").Value = "  -4.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.57"
$ws.Range("E39").Value = "  -6.34%  "

# Row 40
$ws.Range("E40").Value = "  -3.51%  "

# Row 41
$ws.Range("D41").Value = "2.413.99"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.36"
$ws.Range("E42").Value = "  -2.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.31"
$ws.Range("E43").Value = "  -5.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.668"
$ws.Range("E44").Value = "  -2.82%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0596"
$ws.Range("E45").Value = "  -2.42%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.11"
$ws.Range("E46").Value = "  -6.62%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0249"
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.997"
$ws.Range("E48").Value = "  +0.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.96"
$ws.Range("E49").Value = "  -5.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "270.53"
$ws.Range("E50").Value = "  -6.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  -2.47%  "
